$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Existing "2017.12.12" sheet (sheet19.xml): update sprint numbers and add
#    a row for the new "Filter tutors in comment page by subject" task.
# ---------------------------------------------------------------------------
$ws1212 = $wb.Worksheets.Item("2017.12.12")

# Row 3 ("Tutoring request"): estimation/worked numbers change from 6/6/0 to 5/5/5
$ws1212.Cells.Item(3, 3).Value = 5   # C3
$ws1212.Cells.Item(3, 4).Value = 5   # D3
$ws1212.Cells.Item(3, 5).Value = 5   # E3

# New row 4: "Filter tutors in  comment page by subject"
$ws1212.Cells.Item(4, 2).Value = "Filter tutors in  comment page by subject"
$ws1212.Cells.Item(4, 3).Value = 1
$ws1212.Cells.Item(4, 4).Value = 1
$ws1212.Cells.Item(4, 5).Value = 1

# Row 8 ("Images for tutors"): worked goes from 0 to 4 (fully worked)
$ws1212.Cells.Item(8, 5).Value = 4   # E8

# Row 9 ("Design"): worked goes from 0 to 2 (fully worked)
$ws1212.Cells.Item(9, 5).Value = 2   # E9

# Selection on this sheet moves to the whole table range
$ws1212.Range("A1:F17").Select()

# ---------------------------------------------------------------------------
# 2) Brand new "2017.12.19" sprint sheet, added right after "2017.12.12".
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws1219 = $wb.Worksheets.Add($null, $lastSheet)
$ws1219.Name = "2017.12.19"
$ws1219.Columns.Item(2).ColumnWidth = 36

# Header row
$ws1219.Cells.Item(1, 1).Value = "User Story"
$ws1219.Cells.Item(1, 2).Value = "Task"
$ws1219.Cells.Item(1, 3).Value = "Initial Estimation"
$ws1219.Cells.Item(1, 4).Value = "Current Estimation"
$ws1219.Cells.Item(1, 5).Value = "Worked"
$ws1219.Cells.Item(1, 6).Value = "Remain"

# Row 2 + 3: "All" / "Filter bug"
$ws1219.Cells.Item(2, 1).Value = "All"
$ws1219.Cells.Item(3, 2).Value = "Filter bug"
$ws1219.Cells.Item(3, 3).Value = 1
$ws1219.Cells.Item(3, 4).Value = 1
$ws1219.Cells.Item(3, 5).Value = 0
$ws1219.Range("F3").Formula = "=D3-E3"

# Row 4: "Filter tutors in  comment page by subject"
$ws1219.Cells.Item(4, 2).Value = "Filter tutors in  comment page by subject"
$ws1219.Cells.Item(4, 3).Value = 1
$ws1219.Cells.Item(4, 4).Value = 1
$ws1219.Cells.Item(4, 5).Value = 0

# Row 8: "All" / "Admin mail"
$ws1219.Cells.Item(8, 1).Value = "All"
$ws1219.Cells.Item(8, 2).Value = "Admin mail"
$ws1219.Cells.Item(8, 3).Value = 2
$ws1219.Cells.Item(8, 4).Value = 2
$ws1219.Cells.Item(8, 5).Value = 0
$ws1219.Range("F8").Formula = "=D8-E8"

# Row 12: totals
$ws1219.Range("C12").Formula = "=SUM(C2:C11)"
$ws1219.Range("D12").Formula = "=SUM(D2:D11)"
$ws1219.Range("E12").Formula = "=SUM(E2:E11)"
$ws1219.Range("F12").Formula = "=SUM(F2:F11)"

# Row 14-17: planned hours per person
$ws1219.Cells.Item(14, 1).Value = "Name"
$ws1219.Cells.Item(14, 2).Value = "Planned hours"
$ws1219.Cells.Item(15, 1).Value = "Eva"
$ws1219.Cells.Item(15, 2).Value = 2
$ws1219.Cells.Item(16, 1).Value = "Danijal"
$ws1219.Cells.Item(16, 2).Value = 2
$ws1219.Range("B17").Formula = "=SUM(B15:B16)"

$ws1219.Range("G25").Select()

# ---------------------------------------------------------------------------
# 3) Summary sheet: add sprint row for 2017.12.19.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("A22:B22").Copy()
$wsSummary.Range("A23:B23").PasteSpecial(-4122)
$wsSummary.Cells.Item(23, 1).Value = 43088
$wsSummary.Cells.Item(23, 2).Value = "Admin mail, filter"

# Keep "Summary" as the active/selected sheet, matching the original file,
# with the new row's second cell selected (mirrors the prior B22 selection).
$wsSummary.Activate()
$wsSummary.Range("B23").Select()
